$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.779.72'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '3.694.88'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '673.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.81%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.500'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.147'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.82%  '
$ws.Range("E10").Value = '  +2.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.445'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000235'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.74%  '
$ws.Range("D14").Value = '3.682.87'
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").Value = '69.719.72'
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("E16").Value = '  +2.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '16.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '473.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.652'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("D23").Value = '3.842.59'
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("E24").Value = '  +7.54%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.57%  '
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.169'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.35%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("D35").Value = '3.690.26'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.29%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.31'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.79%  '
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0912'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '176.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.938'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '46.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.17%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.48%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000276'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.60%  '
$ws.Range("E49").Value = '  +1.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.92'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.29%  '
$ws.Range("E51").Value = '  +0.61%  '
